$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic_Plotting")

# Add the new "Free_Y" column (F) header and per-parameter T/F values.
$ws.Range("F1").Value = "Free_Y"
$ws.Range("F2").Value  = "T"   # temp
$ws.Range("F3").Value  = "T"   # sal
$ws.Range("F4").Value  = "F"   # do_mgl
$ws.Range("F5").Value  = "T"   # ph
$ws.Range("F6").Value  = "T"   # turb
$ws.Range("F7").Value  = "F"   # depth
$ws.Range("F8").Value  = "T"   # atemp
$ws.Range("F9").Value  = "T"   # rh
$ws.Range("F10").Value = "T"   # bp
$ws.Range("F11").Value = "F"   # wspd
$ws.Range("F12").Value = "F"   # maxwspd
$ws.Range("F13").Value = "F"   # totprcp
$ws.Range("F14").Value = "F"   # totpar
$ws.Range("F15").Value = "F"   # po4f
$ws.Range("F16").Value = "F"   # nh4f
$ws.Range("F17").Value = "F"   # no2f
$ws.Range("F18").Value = "F"   # no3f
$ws.Range("F19").Value = "F"   # no23f
$ws.Range("F20").Value = "F"   # chla_n
$ws.Range("F21").Value = "F"   # din
$ws.Range("F22").Value = "F"   # dip

# Narrow column F now that it only holds short T/F flags (was a wide
# bestFit text column, now a short T/F flag column).
$ws.Columns.Item(6).ColumnWidth = 8.65

# Make Basic_Plotting the active sheet/tab with F1:F22 selected (matches
# the author's final UI state captured in the workbook).
$ws.Activate()
$ws.Range("F1:F22").Select()
